$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("G11").Value = 2.45
$ws.Range("I11").Value = 3.2
$ws.Range("U11").Value = 11
$ws.Range("AF11").Value = 12

# Row 12
$ws.Range("G12").Value = 1.83
$ws.Range("H12").Value = 3.25
$ws.Range("I12").Value = 4.2
$ws.Range("L12").Value = 1.36
$ws.Range("M12").Value = 2.65
$ws.Range("N12").Value = 2.05
$ws.Range("O12").Value = 1.6
$ws.Range("P12").Value = 1.47
$ws.Range("Q12").Value = 2.35
$ws.Range("R12").Value = 1.91
$ws.Range("S12").Value = 1.7
$ws.Range("T12").Value = 5.9
$ws.Range("U12").Value = 7.8
$ws.Range("V12").Value = 8.5
$ws.Range("W12").Value = 15
$ws.Range("X12").Value = 16.5
$ws.Range("Y12").Value = 32
$ws.Range("Z12").Value = 8
$ws.Range("AA12").Value = 6.4
$ws.Range("AB12").Value = 17
$ws.Range("AC12").Value = 90
$ws.Range("AD12").Value = 10.5
$ws.Range("AE12").Value = 23
$ws.Range("AF12").Value = 14
$ws.Range("AG12").Value = 70
$ws.Range("AH12").Value = 45
$ws.Range("AI12").Value = 55
$ws.Range("AJ12").Value = 900

# Row 18
$ws.Range("J18").Value = 1.05
$ws.Range("K18").Value = 11

# Row 19
$ws.Range("N19").Value = 1.95
$ws.Range("O19").Value = 1.85

# Row 21
$ws.Range("I21").Value = 3.7
$ws.Range("K21").Value = 6.4
$ws.Range("U21").Value = 9.25
$ws.Range("Z21").Value = 6.4
$ws.Range("AI21").Value = 50
$ws.Range("AJ21").Value = 700

# Row 22
$ws.Range("I22").Value = 4.35
$ws.Range("J22").Value = 1.06
$ws.Range("K22").Value = 7.5
$ws.Range("L22").Value = 1.27
$ws.Range("M22").Value = 3.4
$ws.Range("N22").Value = 1.8
$ws.Range("O22").Value = 1.9
$ws.Range("R22").Value = 1.72
$ws.Range("S22").Value = 2
$ws.Range("T22").Value = 7.1
$ws.Range("U22").Value = 8.5
$ws.Range("V22").Value = 8
$ws.Range("W22").Value = 14.5
$ws.Range("X22").Value = 14
$ws.Range("Y22").Value = 24
$ws.Range("Z22").Value = 7.5
$ws.Range("AC22").Value = 60
$ws.Range("AD22").Value = 13.5
$ws.Range("AE22").Value = 27
$ws.Range("AF22").Value = 13.5
$ws.Range("AH22").Value = 37
$ws.Range("AI22").Value = 37

# Row 27
$ws.Range("G27").Value = 2.55
$ws.Range("I27").Value = 2.82
$ws.Range("N27").Value = 2.7
$ws.Range("P27").Value = 1.6
$ws.Range("R27").Value = 2.25
$ws.Range("T27").Value = 5.6
$ws.Range("U27").Value = 10.5
$ws.Range("V27").Value = 11.25
$ws.Range("W27").Value = 28
$ws.Range("X27").Value = 30
$ws.Range("AD27").Value = 6
$ws.Range("AE27").Value = 12

# Row 28
$ws.Range("G28").Value = 2.15
$ws.Range("H28").Value = 3.2
$ws.Range("I28").Value = 3.25
$ws.Range("U28").Value = 9
$ws.Range("W28").Value = 19.5
$ws.Range("AA28").Value = 6.3
$ws.Range("AD28").Value = 7.8
$ws.Range("AG28").Value = 45

# Row 31
$ws.Range("G31").Value = 4
$ws.Range("I31").Value = 2
$ws.Range("L31").Value = 1.29
$ws.Range("M31").Value = 3.5
$ws.Range("N31").Value = 1.95
$ws.Range("O31").Value = 1.85
$ws.Range("P31").Value = 1.4
$ws.Range("Q31").Value = 2.75
$ws.Range("R31").Value = 1.8
$ws.Range("S31").Value = 1.95
$ws.Range("T31").Value = 11
$ws.Range("U31").Value = 21
$ws.Range("X31").Value = 34
$ws.Range("Y31").Value = 41
$ws.Range("Z31").Value = 9.5
$ws.Range("AB31").Value = 15
$ws.Range("AC31").Value = 51
$ws.Range("AD31").Value = 7.5
$ws.Range("AE31").Value = 9.5
$ws.Range("AG31").Value = 17
$ws.Range("AJ31").Value = 251

# Row 32
$ws.Range("G32").Value = 2.02
$ws.Range("H32").Value = 3.2
$ws.Range("I32").Value = 3.45
$ws.Range("J32").Value = 1.08
$ws.Range("K32").Value = 6.6
$ws.Range("L32").Value = 1.37
$ws.Range("M32").Value = 2.87
$ws.Range("N32").Value = 2.1
$ws.Range("O32").Value = 1.65
$ws.Range("P32").Value = 1.47
$ws.Range("Q32").Value = 2.5
$ws.Range("R32").Value = 1.88
$ws.Range("S32").Value = 1.82
$ws.Range("T32").Value = 6.6
$ws.Range("U32").Value = 9
$ws.Range("V32").Value = 8.75
$ws.Range("W32").Value = 18
$ws.Range("X32").Value = 17.5
$ws.Range("Y32").Value = 32
$ws.Range("Z32").Value = 6.6
$ws.Range("AA32").Value = 6.3
$ws.Range("AB32").Value = 16
$ws.Range("AC32").Value = 80
$ws.Range("AD32").Value = 9.25
$ws.Range("AE32").Value = 17.5
$ws.Range("AF32").Value = 12
$ws.Range("AG32").Value = 50
$ws.Range("AH32").Value = 35
$ws.Range("AI32").Value = 45
$ws.Range("AJ32").Value = 700
